$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (existing)
$ws.Range("A29").Value = 2920.29
$ws.Range("B29").Value = 3784.65
$ws.Range("C29").Value = 2821.92
$ws.Range("D29").Value = 3757.92
$ws.Range("E29").Value = 2059.11
$ws.Range("F29").Value = 2668.57
$ws.Range("G29").Value = 1989.75
$ws.Range("H29").Value = 2649.73

# Row 32 (existing)
$ws.Range("A32").Value = 2931.23
$ws.Range("B32").Value = 3799.09
$ws.Range("C32").Value = 2834.01
$ws.Range("D32").Value = 3775.39
$ws.Range("E32").Value = 2088.18
$ws.Range("F32").Value = 2707.75
$ws.Range("G32").Value = 2022.05
$ws.Range("H32").Value = 2697.23

# Row 33 (existing)
$ws.Range("A33").Value = 20379.81
$ws.Range("B33").Value = 22549.47
$ws.Range("C33").Value = 20136.76
$ws.Range("D33").Value = 22490.22
$ws.Range("E33").Value = 21620.83
$ws.Range("F33").Value = 23169.76
$ws.Range("G33").Value = 21455.51
$ws.Range("H33").Value = 23143.46
$ws.Range("I33").Value = "p_max_ver"

# Row 38 (existing)
$ws.Range("A38").Value = 20098.83
$ws.Range("B38").Value = 22551.64
$ws.Range("C38").Value = 19926.2
$ws.Range("D38").Value = 22546.78
$ws.Range("E38").Value = 18743.87
$ws.Range("F38").Value = 21197
$ws.Range("G38").Value = 18761.97
$ws.Range("H38").Value = 21454.04

# Row 43 (existing)
$ws.Range("A43").Value = 19855.52
$ws.Range("B43").Value = 22187.94
$ws.Range("C43").Value = 19637.4
$ws.Range("D43").Value = 22194.6
$ws.Range("E43").Value = 26085.51
$ws.Range("F43").Value = 25595.05
$ws.Range("G43").Value = 25477.87
$ws.Range("H43").Value = 24996.05

# Row 48 (existing)
$ws.Range("I48").Value = "номер оси n"

# Row 49 (existing)
$ws.Range("A49").Value = 2009.34
$ws.Range("B49").Value = 1969.08
$ws.Range("C49").Value = 2063.7
$ws.Range("D49").Value = 2001.52
$ws.Range("E49").Value = 1768.94
$ws.Range("F49").Value = 1747.15
$ws.Range("G49").Value = 2001.56
$ws.Range("H49").Value = 1961.79

# Row 50 (existing)
$ws.Range("A50").Value = 11.86
$ws.Range("B50").Value = 15.18
$ws.Range("C50").Value = 11.17
$ws.Range("D50").Value = 14.73
$ws.Range("E50").Value = 15.59
$ws.Range("F50").Value = 17.51
$ws.Range("G50").Value = 14.5
$ws.Range("H50").Value = 16.59

# Row 51 (existing)
$ws.Range("A51").Value = 2.54
$ws.Range("B51").Value = 3.24
$ws.Range("C51").Value = 2.39
$ws.Range("D51").Value = 3.15
$ws.Range("E51").Value = 3.33
$ws.Range("F51").Value = 3.74
$ws.Range("G51").Value = 3.1
$ws.Range("H51").Value = 3.55

# Row 53 (new)
$ws.Range("A53").Value = 0.74428
$ws.Range("B53").Value = 0.68623
$ws.Range("C53").Value = 0.73355
$ws.Range("D53").Value = 0.66465
$ws.Range("E53").Value = 0.74428
$ws.Range("F53").Value = 0.68623
$ws.Range("G53").Value = 0.73355
$ws.Range("H53").Value = 0.66465
$ws.Range("I53").Value = "Тета от длинны шпалы"

# Row 55 (existing)
$ws.Range("A55").Value = 1.26
$ws.Range("B55").Value = 1.15
$ws.Range("C55").Value = 1.29
$ws.Range("D55").Value = 1.16
$ws.Range("E55").Value = 1.14
$ws.Range("F55").Value = 1.08
$ws.Range("G55").Value = 1.17
$ws.Range("H55").Value = 1.11

# Row 56 (existing)
$ws.Range("A56").Value = 0.527
$ws.Range("B56").Value = 0.671
$ws.Range("C56").Value = 0.496
$ws.Range("D56").Value = 0.652
$ws.Range("E56").Value = 0.689
$ws.Range("F56").Value = 0.774
$ws.Range("G56").Value = 0.642
$ws.Range("H56").Value = 0.735
$ws.Range("I56").Value = "sigma_h2()"

# Row 57 (existing)
$ws.Range("A57").Value = 0.248
$ws.Range("B57").Value = 0.248
$ws.Range("C57").Value = 0.248
$ws.Range("D57").Value = 0.248
$ws.Range("E57").Value = 0.248
$ws.Range("F57").Value = 0.248
$ws.Range("G57").Value = 0.248
$ws.Range("H57").Value = 0.248
$ws.Range("I57").Value = "А, коэффициент расстояния между шпал"

# Row 58 (existing)
$ws.Range("A58").Value = 0.206
$ws.Range("B58").Value = 0.206
$ws.Range("C58").Value = 0.206
$ws.Range("D58").Value = 0.206
$ws.Range("E58").Value = 0.206
$ws.Range("F58").Value = 0.206
$ws.Range("G58").Value = 0.206
$ws.Range("H58").Value = 0.206
$ws.Range("I58").Value = "C1"

# Row 59 (existing)
$ws.Range("I59").Value = "l_i[0]"

# Row 60 (existing)
$ws.Range("I60").Value = "l_i[1]"

# Row 61 (existing)
$ws.Range("I61").Value = "l_i[2]"

# Row 62 (new)
$ws.Range("A62").Value = 0.9
$ws.Range("B62").Value = 0.9
$ws.Range("C62").Value = 0.9
$ws.Range("D62").Value = 0.9
$ws.Range("E62").Value = 0.9
$ws.Range("F62").Value = 0.9
$ws.Range("G62").Value = 0.9
$ws.Range("H62").Value = 0.9
$ws.Range("I62").Value = "[бз_Вагон]"

# Row 64 (existing)
$ws.Range("A64").Value = 0.1
$ws.Range("B64").Value = 0.1
$ws.Range("C64").Value = 0.1
$ws.Range("D64").Value = 0.1
$ws.Range("E64").Value = 0.1
$ws.Range("F64").Value = 0.1
$ws.Range("G64").Value = 0.1
$ws.Range("H64").Value = 0.1
$ws.Range("I64").Value = "C2"

# Row 67 (existing)
$ws.Range("A67").Value = 55
$ws.Range("B67").Value = 57
$ws.Range("C67").Value = 52
$ws.Range("D67").Value = 55
$ws.Range("E67").Value = 67
$ws.Range("F67").Value = 69
$ws.Range("G67").Value = 55
$ws.Range("H67").Value = 58

# Row 81 (existing)
$ws.Range("I81").Value = "xn"

# Row 82 (existing)
$ws.Range("A82").Value = -0.02405
$ws.Range("B82").Value = -0.01006
$ws.Range("C82").Value = -0.02108
$ws.Range("D82").Value = -0.00636
$ws.Range("E82").Value = 0.27492
$ws.Range("F82").Value = 0.14787
$ws.Range("G82").Value = 0.24803
$ws.Range("H82").Value = 0.11076
$ws.Range("I82").Value = "сигма тета 1 шпала"

# Row 95 (existing)
$ws.Range("A95").Value = -0.03714
$ws.Range("B95").Value = -0.04301
$ws.Range("C95").Value = -0.03974
$ws.Range("D95").Value = -0.04142
$ws.Range("E95").Value = 0.55906
$ws.Range("F95").Value = 0.47506
$ws.Range("G95").Value = 0.54038
$ws.Range("H95").Value = 0.45179
$ws.Range("I95").Value = "сигма тета 3 шпала"

# Row 100 (existing)
$ws.Range("A100").Value = 14854.39
$ws.Range("B100").Value = 15342.82
$ws.Range("C100").Value = 14496.19
$ws.Range("D100").Value = 14865.12
$ws.Range("E100").Value = 20600.74
$ws.Range("F100").Value = 18324.91
$ws.Range("G100").Value = 19806.48
$ws.Range("H100").Value = 17198.81
$ws.Range("I100").Value = "P_II_ekvONE"

# Row 101 (existing)
$ws.Range("A101").Value = 14683.54
$ws.Range("B101").Value = 14912.77
$ws.Range("C101").Value = 14252.64
$ws.Range("D101").Value = 14407.52
$ws.Range("E101").Value = 25260.75
$ws.Range("F101").Value = 23690.95
$ws.Range("G101").Value = 24601.13
$ws.Range("H101").Value = 22791.83
$ws.Range("I101").Value = "P_II_ekvThree"

# Row 102 (existing)
$ws.Range("A102").Value = 2.63
$ws.Range("B102").Value = 2.63
$ws.Range("C102").Value = 2.63
$ws.Range("D102").Value = 2.63
$ws.Range("E102").Value = 2.63
$ws.Range("F102").Value = 2.63
$ws.Range("G102").Value = 2.63
$ws.Range("H102").Value = 2.63
$ws.Range("I102").Value = "sigma_b1"

# Row 103 (existing)
$ws.Range("A103").Value = 0.1631
$ws.Range("B103").Value = 0.1631
$ws.Range("C103").Value = 0.1631
$ws.Range("D103").Value = 0.1631
$ws.Range("E103").Value = 0.1631
$ws.Range("F103").Value = 0.1631
$ws.Range("G103").Value = 0.1631
$ws.Range("H103").Value = 0.1631
$ws.Range("I103").Value = "sigma_h1"

# Row 104 (existing)
$ws.Range("A104").Value = 3.225
$ws.Range("B104").Value = 3.225
$ws.Range("C104").Value = 3.225
$ws.Range("D104").Value = 3.225
$ws.Range("E104").Value = 3.225
$ws.Range("F104").Value = 3.225
$ws.Range("G104").Value = 3.225
$ws.Range("H104").Value = 3.225
$ws.Range("I104").Value = "sigma_b3"

# Row 105 (existing)
$ws.Range("A105").Value = 0.2
$ws.Range("B105").Value = 0.2
$ws.Range("C105").Value = 0.2
$ws.Range("D105").Value = 0.2
$ws.Range("E105").Value = 0.2
$ws.Range("F105").Value = 0.2
$ws.Range("G105").Value = 0.2
$ws.Range("H105").Value = 0.2
$ws.Range("I105").Value = "sigma_h3"

# Row 106 (existing)
$ws.Range("A106").Value = 1.052
$ws.Range("B106").Value = 1.052
$ws.Range("C106").Value = 1.052
$ws.Range("D106").Value = 1.052
$ws.Range("E106").Value = 1.052
$ws.Range("F106").Value = 1.052
$ws.Range("G106").Value = 1.052
$ws.Range("H106").Value = 1.052
$ws.Range("I106").Value = "∑_h"

# Row 113 (existing)
$ws.Range("A113").Value = 57
$ws.Range("B113").Value = 57
$ws.Range("C113").Value = 57
$ws.Range("D113").Value = 57
$ws.Range("E113").Value = 57
$ws.Range("F113").Value = 57
$ws.Range("G113").Value = 57
$ws.Range("H113").Value = 57

# Row 114 (existing)
$ws.Range("A114").Value = 3996.2
$ws.Range("B114").Value = 3996.2
$ws.Range("C114").Value = 3996.2
$ws.Range("D114").Value = 3996.2
$ws.Range("E114").Value = 3996.2
$ws.Range("F114").Value = 3996.2
$ws.Range("G114").Value = 3996.2
$ws.Range("H114").Value = 3996.2

# Row 115 (existing)
$ws.Range("A115").Value = 55
$ws.Range("B115").Value = 55
$ws.Range("C115").Value = 55
$ws.Range("D115").Value = 55
$ws.Range("E115").Value = 55
$ws.Range("F115").Value = 55
$ws.Range("G115").Value = 55
$ws.Range("H115").Value = 55

# Row 116 (existing)
$ws.Range("A116").Value = 3998.14
$ws.Range("B116").Value = 3998.14
$ws.Range("C116").Value = 3998.14
$ws.Range("D116").Value = 3998.14
$ws.Range("E116").Value = 3998.14
$ws.Range("F116").Value = 3998.14
$ws.Range("G116").Value = 3998.14
$ws.Range("H116").Value = 3998.14

# Row 117 (existing)
$ws.Range("A117").Value = 1436.4
$ws.Range("B117").Value = 1436.4
$ws.Range("C117").Value = 1436.4
$ws.Range("D117").Value = 1436.4
$ws.Range("E117").Value = 1436.4
$ws.Range("F117").Value = 1436.4
$ws.Range("G117").Value = 1436.4
$ws.Range("H117").Value = 1436.4

# Row 118 (existing)
$ws.Range("A118").Value = 1386
$ws.Range("B118").Value = 1386
$ws.Range("C118").Value = 1386
$ws.Range("D118").Value = 1386
$ws.Range("E118").Value = 1386
$ws.Range("F118").Value = 1386
$ws.Range("G118").Value = 1386
$ws.Range("H118").Value = 1386

# Row 121 (existing)
$ws.Range("A121").Value = 1969.08
$ws.Range("B121").Value = 1969.08
$ws.Range("C121").Value = 1969.08
$ws.Range("D121").Value = 1969.08
$ws.Range("E121").Value = 1969.08
$ws.Range("F121").Value = 1969.08
$ws.Range("G121").Value = 1969.08
$ws.Range("H121").Value = 1969.08

# Row 122 (existing)
$ws.Range("A122").Value = 2009.34
$ws.Range("B122").Value = 2009.34
$ws.Range("C122").Value = 2009.34
$ws.Range("D122").Value = 2009.34
$ws.Range("E122").Value = 2009.34
$ws.Range("F122").Value = 2009.34
$ws.Range("G122").Value = 2009.34
$ws.Range("H122").Value = 2009.34

# Row 124 (existing)
$ws.Range("A124").Value = 2729.44
$ws.Range("B124").Value = 2729.44
$ws.Range("C124").Value = 2729.44
$ws.Range("D124").Value = 2729.44
$ws.Range("E124").Value = 2729.44
$ws.Range("F124").Value = 2729.44
$ws.Range("G124").Value = 2729.44
$ws.Range("H124").Value = 2729.44

# Row 126 (existing)
$ws.Range("A126").Value = 278325.42
$ws.Range("B126").Value = 278325.42
$ws.Range("C126").Value = 278325.42
$ws.Range("D126").Value = 278325.42
$ws.Range("E126").Value = 278325.42
$ws.Range("F126").Value = 278325.42
$ws.Range("G126").Value = 278325.42
$ws.Range("H126").Value = 278325.42

# Row 128 (existing)
$ws.Range("A128").Value = 185550.28
$ws.Range("B128").Value = 185550.28
$ws.Range("C128").Value = 185550.28
$ws.Range("D128").Value = 185550.28
$ws.Range("E128").Value = 185550.28
$ws.Range("F128").Value = 185550.28
$ws.Range("G128").Value = 185550.28
$ws.Range("H128").Value = 185550.28

# Row 131 (existing)
$ws.Range("A131").Value = 65
$ws.Range("B131").Value = 65
$ws.Range("C131").Value = 65
$ws.Range("D131").Value = 65
$ws.Range("E131").Value = 65
$ws.Range("F131").Value = 65
$ws.Range("G131").Value = 65
$ws.Range("H131").Value = 65
$ws.Range("I131").Value = "[∆t_уПрямая]"

# Row 132 (existing)
$ws.Range("A132").Value = 56
$ws.Range("B132").Value = 56
$ws.Range("C132").Value = 56
$ws.Range("D132").Value = 56
$ws.Range("E132").Value = 56
$ws.Range("F132").Value = 56
$ws.Range("G132").Value = 56
$ws.Range("H132").Value = 56
$ws.Range("I132").Value = "[∆t_у_curve]"

# Row 133 (existing)
$ws.Range("A133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("B133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("C133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("D133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("E133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("F133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("G133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("H133").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; η = -0.04017`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; η = 0`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; η = 0"
$ws.Range("I133").Value = "Ekv_gruzi_η"

# Row 134 (existing)
$ws.Range("A134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("B134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("C134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("D134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("E134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("F134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("G134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("H134").Value = "II ось: x = 300 см; kx = 0.01145×300 = 3.44; µ = -0.02153`nIII ось: x = 0+300 см; kx = 0.01145×300 = 3.44; µ = 0.00000`nVI ось: x = 300+0+0 см; kx = 0.01145×300 = 3.44; µ = 0.00000"
$ws.Range("I134").Value = "Ekv_gruzi_µ"

# Row 135 (existing)
$ws.Range("A135").Value = "ηI: x = 185 - 55 см; kx = 0.01145×130 = 1.49; η = 0.24350`nηII: x = 55 см; kx = 0.01145×55 = 0.63; η = 0.74428`nηIII: x = 125+55 см; kx = 0.01145×180 = 2.06; η = 0.05239`nηIV: x = 125+185+55; kx = 0.01145×365 = 4.18; η = -0.02097"

# Row 136 (new)
$ws.Range("A136").Value = "ηI: x = 185 + 55 см; kx = 0.01145×240 = 2.75; η = -0.03459`nηII: x = 55 см; kx = 0.01145×55 = 0.63; η = 0.74428`nηIII: x = 125 - 55 см; kx = 0.01145×70 = 0.80; η = 0.63441`nηIV: x = 125 + 185-55; kx = 0.01145×255 = 3.49; η = -0.04076"
$ws.Range("I136").Value = "Ekv_gruzi_η_shpala_3"

# Row 137 (new)
$ws.Range("A137").Value = 58
$ws.Range("B137").Value = 58
$ws.Range("C137").Value = 58
$ws.Range("D137").Value = 58
$ws.Range("E137").Value = 58
$ws.Range("F137").Value = 58
$ws.Range("G137").Value = 58
$ws.Range("H137").Value = 58
$ws.Range("I137").Value = "t_max_max"

# Row 138 (new)
$ws.Range("A138").Value = -42
$ws.Range("B138").Value = -42
$ws.Range("C138").Value = -42
$ws.Range("D138").Value = -42
$ws.Range("E138").Value = -42
$ws.Range("F138").Value = -42
$ws.Range("G138").Value = -42
$ws.Range("H138").Value = -42
$ws.Range("I138").Value = "t_min_min"

# Row 139 (new)
$ws.Range("A139").Value = 100
$ws.Range("B139").Value = 100
$ws.Range("C139").Value = 100
$ws.Range("D139").Value = 100
$ws.Range("E139").Value = 100
$ws.Range("F139").Value = 100
$ws.Range("G139").Value = 100
$ws.Range("H139").Value = 100
$ws.Range("I139").Value = "Tа"

# Row 140 (new)
$ws.Range("A140").Value = 101.24
$ws.Range("B140").Value = 101.24
$ws.Range("C140").Value = 101.24
$ws.Range("D140").Value = 101.24
$ws.Range("E140").Value = 101.24
$ws.Range("F140").Value = 101.24
$ws.Range("G140").Value = 101.24
$ws.Range("H140").Value = 101.24
$ws.Range("I140").Value = "[T] кривая"

# Row 141 (new)
$ws.Range("A141").Value = 112.15
$ws.Range("B141").Value = 112.15
$ws.Range("C141").Value = 112.15
$ws.Range("D141").Value = 112.15
$ws.Range("E141").Value = 112.15
$ws.Range("F141").Value = 112.15
$ws.Range("G141").Value = 112.15
$ws.Range("H141").Value = 112.15
$ws.Range("I141").Value = "[T] прямая"

# Row 142 (new)
$ws.Range("A142").Value = 15
$ws.Range("B142").Value = 15
$ws.Range("C142").Value = 15
$ws.Range("D142").Value = 15
$ws.Range("E142").Value = 15
$ws.Range("F142").Value = 15
$ws.Range("G142").Value = 15
$ws.Range("H142").Value = 15
$ws.Range("I142").Value = "t_max_zakr"

# Row 143 (new)
$ws.Range("A143").Value = 13
$ws.Range("B143").Value = 13
$ws.Range("C143").Value = 13
$ws.Range("D143").Value = 13
$ws.Range("E143").Value = 13
$ws.Range("F143").Value = 13
$ws.Range("G143").Value = 13
$ws.Range("H143").Value = 13
$ws.Range("I143").Value = "t_max_zakr_curve"

# Row 144 (new)
$ws.Range("A144").Value = -7
$ws.Range("B144").Value = -7
$ws.Range("C144").Value = -7
$ws.Range("D144").Value = -7
$ws.Range("E144").Value = -7
$ws.Range("F144").Value = -7
$ws.Range("G144").Value = -7
$ws.Range("H144").Value = -7
$ws.Range("I144").Value = "t_min_zakr"

# Row 145 (new)
$ws.Range("A145").Value = 2
$ws.Range("B145").Value = 2
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 2
$ws.Range("E145").Value = 2
$ws.Range("F145").Value = 2
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 2
$ws.Range("I145").Value = "t_min_zakr_curve"
